$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = 44959
$ws.Range("A8").NumberFormat = "d-mmm"
$ws.Range("B8").Value = "10.30 - 13.00"
$ws.Range("D8").Value = 2.5
$ws.Range("E8").Value = "Implementing dijkstra algorithm"

$ws.Range("A11").Value = 44967
$ws.Range("A11").NumberFormat = "d-mmm"
$ws.Range("B11").Value = "9.30 - 14.30"
$ws.Range("D11").Value = 5

$ws.Range("A9").Value = 44965
$ws.Range("A9").NumberFormat = "d-mmm"
$ws.Range("E9").Value = "Implemented cost field, heatmap"
$ws.Range("B9").Value = "10.00 - 12.30"
$ws.Range("D9").Value = 2.5

$ws.Range("A10").Value = 44966
$ws.Range("A10").NumberFormat = "d-mmm"
$ws.Range("B10").Value = "9.00 - 12.30"
$ws.Range("D10").Value = 3.5
$ws.Range("E10").Value = "Implemented flow field, working on AI agents"

$ws.Range("E13").Select()
